$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames / additions -------------------------------------------------
$ws.Range("G1").Value = "Height_m"
$ws.Range("H1").Value = "BMI_Baseline"
$ws.Range("I1").Value = "BMI_Post_Int"
$ws.Range("J1").Value = "BMI_Followup"

# --- Fix mis-entered PPID values in rows 41-50 (were all 39) -------------------
for ($r = 41; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- Complete row 51 (new participant) ------------------------------------------
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = 2
$ws.Range("C51").Value = 33
$ws.Range("G51").Value = 1.82

# --- New BMI formula columns -----------------------------------------------------
# Row 2 is its own (non-shared) anchor formula, rows 3-51 fill down as one
# shared-formula block per column (mirrors how Excel itself lays these out).
$ws.Range("H51").Formula = "=(D51/(G51*G51))"

$ws.Range("I2").Formula = "=(E2/(G2*G2))"
$ws.Range("I3:I51").Formula = "=(E3/(G3*G3))"

$ws.Range("J2").Formula = "=(F2/(G2*G2))"
$ws.Range("J3:J51").Formula = "=(F3/(G3*G3))"

# --- Number formatting for the BMI columns ---------------------------------------
$ws.Range("H2:J51").NumberFormat = "0.00"

# --- Match the cursor position left behind by the edit session ------------------
[void]$ws.Range("N50").Select()
